$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# Metadata sheet: Title row (A5="Title"), Description row (A12="Description"), Date row (A8="Date")
$meta.Range("B5").Value = "DMI Code LPP"
$meta.Range("B12").Value = "Extension créée dans ce volet pour représenter le code LPP."
$meta.Range("B8").Value = "2026-02-25T08:15:31+00:00"

# Elements sheet: root Extension element row (row 2) - Short (L2) and Definition (M2)
$elements.Range("L2").Value = "DMI Code LPP"
$elements.Range("M2").Value = "Extension créée dans ce volet pour représenter le code LPP."
